$wb = $excel.ActiveWorkbook

# REPCA1 sheet: insert three new parameter columns (VCFlag, RefFlag, Fflag)
# right after the existing "busf" column (H), shifting Tfltr and everything
# after it three columns to the right.
$ws = $wb.Worksheets.Item("REPCA1")
$ws.Range("I1:K1").EntireColumn.Insert()

$ws.Cells.Item(1, 9).Value = "VCFlag"
$ws.Cells.Item(1, 10).Value = "RefFlag"
$ws.Cells.Item(1, 11).Value = "Fflag"

$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0

# Restore default (top-left) selection/active-tab state: the workbook was
# resaved with the first sheet ("Bus") active and no sheet carrying a
# custom selection anymore.
$wsBus = $wb.Worksheets.Item("Bus")
[void]$wsBus.Activate()
[void]$wsBus.Range("A1").Select()

$wsBusFreq = $wb.Worksheets.Item("BusFreq")
[void]$wsBusFreq.Range("A1").Select()

[void]$ws.Range("A1").Select()

[void]$wsBus.Activate()
